$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata updates ---
$ws.Range("A2").Value = "Independent Power Producers by Census Division and State, November 2016 (Continued)"

# --- Data cell updates (RSE table rows 4-65) ---
# Row 4
$ws.Range("E4").Value = 11
$ws.Range("H4").Value = 5

# Row 5
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 5
$ws.Range("H5").Value = 7

# Row 6
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 5

# Row 7
$ws.Range("E7").Value = 11
$ws.Range("F7").Value = 8
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 3

# Row 8
$ws.Range("F8").Value = 13
$ws.Range("H8").Value = 45

# Row 9
$ws.Range("E9").Value = 62
$ws.Range("F9").Value = 17
$ws.Range("I9").Value = 1

# Row 10
$ws.Range("E10").Value = 34
$ws.Range("F10").Value = 34
$ws.Range("I10").Value = 36

# Row 11
$ws.Range("E11").Value = 10
$ws.Range("H11").Value = 4

# Row 12
$ws.Range("E12").Value = 12
$ws.Range("H12").Value = 10

# Row 13
$ws.Range("H13").Value = 6

# Row 14
$ws.Range("E14").Value = 33

# Row 15
$ws.Range("E15").Value = 15
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 0.22

# Row 16
$ws.Range("E16").Value = 35

# Row 17
$ws.Range("E17").Value = 18

# Row 18
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 25

# Row 19
$ws.Range("E19").Value = 35
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 1819
$ws.Range("I19").Value = 0.4

# Row 20
$ws.Range("E20").Value = 201
$ws.Range("F20").Value = 8

# Row 21
$ws.Range("E21").Value = 45
$ws.Range("F21").Value = 0.49
$ws.Range("H21").Value = 29

# Row 23
$ws.Range("E23").Value = 204
$ws.Range("F23").Value = 0.3
$ws.Range("I23").Value = 0.39

# Row 24
$ws.Range("E24").Value = 101
$ws.Range("H24").Value = 29

# Row 25
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 2
$ws.Range("I25").Value = 3

# Row 26
$ws.Range("E26").Value = 102
$ws.Range("F26").Value = 0.41
$ws.Range("I26").Value = 0.41

# Row 27
$ws.Range("F27").Value = 0.34
$ws.Range("I27").Value = 0.34

# Row 29
$ws.Range("E29").Value = 4

# Row 30
$ws.Range("E30").Value = 39
$ws.Range("F30").Value = 32
$ws.Range("I30").Value = 8

# Row 32
$ws.Range("E32").Value = 31
$ws.Range("H32").Value = 5
$ws.Range("I32").Value = 6

# Row 33
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 6
$ws.Range("I33").Value = 4

# Row 34
$ws.Range("E34").Value = 20

# Row 35
$ws.Range("I35").Value = 4

# Row 36
$ws.Range("E36").Value = 121
$ws.Range("F36").Value = 40
$ws.Range("H36").Value = 174
$ws.Range("I36").Value = 4

# Row 37
$ws.Range("F37").Value = 7

# Row 38
$ws.Range("F38").Value = 0.47

# Row 39
$ws.Range("E39").Value = 9
$ws.Range("I39").Value = 0.13

# Row 40
$ws.Range("F40").Value = 5
$ws.Range("I40").Value = 0.06

# Row 41
$ws.Range("F41").Value = 133
$ws.Range("I41").Value = 12

# Row 42
$ws.Range("F42").Value = 100
$ws.Range("I42").Value = 0.36

# Row 43
$ws.Range("E43").Value = 29
$ws.Range("F43").Value = 20
$ws.Range("I43").Value = 20

# Row 44
$ws.Range("E44").Value = 6
$ws.Range("F44").Value = 0.33
$ws.Range("H44").Value = 45
$ws.Range("I44").Value = 0.24

# Row 45
$ws.Range("E45").Value = 59

# Row 46
$ws.Range("F46").Value = 37
$ws.Range("I46").Value = 0.49

# Row 47
$ws.Range("F47").Value = 0.21
$ws.Range("I47").Value = 0.4

# Row 48
$ws.Range("E48").Value = 6
$ws.Range("F48").Value = 0.43
$ws.Range("H48").Value = 99
$ws.Range("I48").Value = 0.28

# Row 49
$ws.Range("H49").Value = 4
$ws.Range("I49").Value = 2

# Row 50
$ws.Range("E50").Value = 2
$ws.Range("F50").Value = 2
$ws.Range("I50").Value = 1

# Row 51
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = 0.41
$ws.Range("H51").Value = 102

# Row 52
$ws.Range("C52").Value = 48
$ws.Range("F52").Value = 4
$ws.Range("I52").Value = 6

# Row 53
$ws.Range("F53").Value = 2

# Row 54
$ws.Range("E54").Value = 2
$ws.Range("F54").Value = 4

# Row 55
$ws.Range("C55").Value = 122
$ws.Range("E55").Value = 5

# Row 56
$ws.Range("H56").Value = 228
$ws.Range("I56").Value = 13

# Row 57
$ws.Range("I57").Value = 20

# Row 58
$ws.Range("H58").Value = 15

# Row 59
$ws.Range("F59").Value = 2
$ws.Range("H59").Value = 16

# Row 60
$ws.Range("E60").Value = 24
$ws.Range("F60").Value = 2
$ws.Range("H60").Value = 49

# Row 61
$ws.Range("H61").Value = 43
$ws.Range("I61").Value = 2

# Row 62
$ws.Range("E62").Value = 39
$ws.Range("I62").Value = 6

# Row 63
$ws.Range("F63").Value = 41
$ws.Range("I63").Value = 46

# Row 64
$ws.Range("E64").Value = 39
$ws.Range("I64").Value = 6

# Row 65
$ws.Range("F65").Value = 0.46
$ws.Range("H65").Value = 3
$ws.Range("I65").Value = 0.23
